$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.680.72"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "1.694.86"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.68%  "
$ws.Range("D5").Value = "'241.91"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "'0.4909"
$ws.Range("E7").Value = "  -4.92%  "
$ws.Range("D8").Value = "'0.2661"
$ws.Range("E8").Value = "  -3.28%  "
$ws.Range("D9").Value = "'0.06049"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "1.720.58"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "'0.07173"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.6344"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'14.74"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "'4.659"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").Value = "'75.10"
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "'1.009"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "25.713.48"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "'11.65"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "'0.000006680"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "1.944.44"
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.517"
$ws.Range("E22").Value = "  +5.51%  "
$ws.Range("D23").Value = "'8.687"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "'5.348"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").Value = "'133.43"
$ws.Range("E25").Value = "  -4.08%  "
$ws.Range("D26").Value = "'14.97"
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").Value = "'1.405"
$ws.Range("E27").Value = "  -7.24%  "
$ws.Range("D28").Value = "'1.742"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").Value = "'103.85"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").Value = "'3.850"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").Value = "'0.08024"
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("D32").Value = "'3.593"
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").Value = "'0.04646"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("D34").Value = "'2.676"
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("D35").Value = "'0.9712"
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("D36").Value = "'0.5962"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("D37").Value = "'2.686"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").Value = "'0.01577"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "'0.8520"
$ws.Range("E39").Value = "  +15.18%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "'1.006"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'1.900"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").Value = "'99.20"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").Value = "'0.3796"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "'4.929"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").Value = "'0.1157"
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("D46").Value = "'6.202"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05251"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'54.47"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").Value = "'29.99"
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("D50").Value = "'7.526"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").Value = "'0.3391"
$ws.Range("E51").Value = "  -0.09%  "
